# Update bacnetDeviceId column (C) values from 0 to 1 for rows 2-4,
# and move the active cell selection to C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1

$ws.Range("C5").Select()
